$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Last updated: 2025-07-15 10:52:33'
$ws.Range("A3").Value = '3T/PO252272'
$ws.Range("B3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("A4").Value = '4010016033'
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("A5").Value = '4020007186'
$ws.Range("B5").Value = 11
$ws.Range("G5").Value = 11
$ws.Range("A6").Value = '4516260169'
$ws.Range("B6").Value = 7
$ws.Range("G6").Value = 7
$ws.Range("A7").Value = '4516351202_AIZU'
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = $null
$ws.Range("G7").Value = 12
$ws.Range("I7").Value = -1
$ws.Range("A8").Value = '4516351202_ARD'
$ws.Range("B8").Value = 25
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 8
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = 16
$ws.Range("I8").Value = -1
$ws.Range("A9").Value = '4516351202_ASEWH'
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 'Completed PO'
$ws.Range("G9").Value = 14
$ws.Range("I9").Value = 0
$ws.Range("A10").Value = '4516351202_ATK'
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 'Completed PO'
$ws.Range("G10").Value = 19
$ws.Range("I10").Value = 0
$ws.Range("A11").Value = '4516351202_DMOS5'
$ws.Range("B11").Value = 2
$ws.Range("G11").Value = 2
$ws.Range("A12").Value = '4516351202_HNT'
$ws.Range("B12").Value = 12
$ws.Range("G12").Value = 12
$ws.Range("A13").Value = '4516351202_LFAB'
$ws.Range("B13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("G13").Value = 0
$ws.Range("A14").Value = '4516351202_MIHO'
$ws.Range("B14").Value = 31
$ws.Range("G14").Value = 31
$ws.Range("A15").Value = '4516351202_NFME'
$ws.Range("B15").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("A16").Value = '4516351202_SCT'
$ws.Range("B16").Value = 29
$ws.Range("G16").Value = 29
$ws.Range("A17").Value = '4516351202_SFAB'
$ws.Range("B17").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("G17").Value = 0
$ws.Range("B24").Value = 4
$ws.Range("D24").Value = 2
